$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 - H (meta_avg) and I (meta_std)
$ws.Range("H4").Value = 756.0699018214594
$ws.Range("I4").Value = 1187.801294409434

# Row 5 - H (meta_avg), I (meta_std), J (meta_min)
$ws.Range("H5").Value = 1196.161078194156
$ws.Range("I5").Value = 1922.757443193122
$ws.Range("J5").Value = 0
